$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Assigns a value to a cell while preserving text ("inline string") semantics.
# Values that look like plain numbers would otherwise be auto-converted by Excel
# into numeric cells, so those are written with a leading apostrophe (the normal
# Excel convention for forcing text entry); the stored cell value has no apostrophe.
function Set-TextValue($range, $value) {
    if ($value -match '^[+-]?\d+(\.\d+)?$') {
        $range.Value = "'" + $value
    } else {
        $range.Value = $value
    }
}

Set-TextValue $ws.Range("D2") '62.933.14'
Set-TextValue $ws.Range("E2") '  -1.21%  '
Set-TextValue $ws.Range("D3") '3.152.36'
Set-TextValue $ws.Range("E3") '  +0.89%  '
Set-TextValue $ws.Range("E4") '  +0.01%  '
Set-TextValue $ws.Range("D5") '588.17'
Set-TextValue $ws.Range("E5") '  -1.94%  '
Set-TextValue $ws.Range("D6") '138.02'
Set-TextValue $ws.Range("E6") '  -3.04%  '
Set-TextValue $ws.Range("E7") '  -0.09%  '
Set-TextValue $ws.Range("D8") '3.152.70'
Set-TextValue $ws.Range("E8") '  +1.06%  '
Set-TextValue $ws.Range("D9") '0.515'
Set-TextValue $ws.Range("E9") '  -1.20%  '
Set-TextValue $ws.Range("E10") '  -1.69%  '
Set-TextValue $ws.Range("E11") '  -0.89%  '
Set-TextValue $ws.Range("E12") '  -1.95%  '
Set-TextValue $ws.Range("D13") '0.0000243'
Set-TextValue $ws.Range("E13") '  -3.88%  '
Set-TextValue $ws.Range("D14") '34.01'
Set-TextValue $ws.Range("E14") '  -2.91%  '
Set-TextValue $ws.Range("D15") '3.670.25'
Set-TextValue $ws.Range("E15") '  +0.94%  '
Set-TextValue $ws.Range("E16") '  +0.87%  '
Set-TextValue $ws.Range("D17") '3.153.30'
Set-TextValue $ws.Range("E17") '  +0.80%  '
Set-TextValue $ws.Range("D18") '62.893.43'
Set-TextValue $ws.Range("E18") '  -1.42%  '
Set-TextValue $ws.Range("E19") '  -2.47%  '
Set-TextValue $ws.Range("D20") '476.22'
Set-TextValue $ws.Range("E20") '  -1.42%  '
Set-TextValue $ws.Range("D21") '13.90'
Set-TextValue $ws.Range("E21") '  -5.45%  '
Set-TextValue $ws.Range("E22") '  -1.00%  '
Set-TextValue $ws.Range("E23") '  +1.49%  '
Set-TextValue $ws.Range("D24") '84.44'
Set-TextValue $ws.Range("E24") '  -2.57%  '
Set-TextValue $ws.Range("D25") '12.95'
Set-TextValue $ws.Range("E25") '  -2.96%  '
Set-TextValue $ws.Range("E26") '  -0.03%  '
Set-TextValue $ws.Range("E27") '  -1.51%  '
Set-TextValue $ws.Range("D28") '6.99'
Set-TextValue $ws.Range("E28") '  +0.72%  '
Set-TextValue $ws.Range("D29") '7.91'
Set-TextValue $ws.Range("E29") '  -3.88%  '
Set-TextValue $ws.Range("D30") '2.07'
Set-TextValue $ws.Range("E30") '  +1.62%  '
Set-TextValue $ws.Range("E31") '  +0.10%  '
Set-TextValue $ws.Range("D32") '26.75'
Set-TextValue $ws.Range("E32") '  -0.94%  '
Set-TextValue $ws.Range("D33") '0.105'
Set-TextValue $ws.Range("E33") '  -4.86%  '
Set-TextValue $ws.Range("D34") '2.51'
Set-TextValue $ws.Range("E34") '  -5.57%  '
Set-TextValue $ws.Range("D35") '1.08'
Set-TextValue $ws.Range("E35") '  -2.35%  '
Set-TextValue $ws.Range("D36") '52.52'
Set-TextValue $ws.Range("E36") '  -0.10%  '
Set-TextValue $ws.Range("E37") '  -3.18%  '
Set-TextValue $ws.Range("D38") '0.0₃0699'
Set-TextValue $ws.Range("E38") '  -5.59%  '
Set-TextValue $ws.Range("D39") '0.0385'
Set-TextValue $ws.Range("E39") '  -2.26%  '
Set-TextValue $ws.Range("D40") '416.30'
Set-TextValue $ws.Range("E40") '  -4.72%  '
Set-TextValue $ws.Range("D41") '2.74'
Set-TextValue $ws.Range("E41") '  -6.61%  '
Set-TextValue $ws.Range("D42") '2.949.50'
Set-TextValue $ws.Range("E42") '  +2.95%  '
Set-TextValue $ws.Range("E43") '  -0.01%  '
Set-TextValue $ws.Range("E44") '  -7.45%  '
Set-TextValue $ws.Range("D45") '0.257'
Set-TextValue $ws.Range("E45") '  -0.50%  '
Set-TextValue $ws.Range("E46") '  +0.04%  '
Set-TextValue $ws.Range("E47") '  -3.32%  '
Set-TextValue $ws.Range("D48") '25.35'
Set-TextValue $ws.Range("E48") '  -1.85%  '
Set-TextValue $ws.Range("E49") '  -0.54%  '
Set-TextValue $ws.Range("D50") '2.24'
Set-TextValue $ws.Range("E50") '  -6.10%  '
Set-TextValue $ws.Range("D51") '120.08'
Set-TextValue $ws.Range("E51") '  -1.73%  '
